# Update the two "square adjacency matrix" equations (Figure 1 and the
# subsequent "after one iteration" figure) with corrected distance values,
# per the iterative-function worked example fix.
$d = $word.ActiveDocument

$mathBlocks = @(
    @{ Index = 1; Xml = '<m:oMathPara><m:oMath><m:d><m:dPr><m:begChr m:val="["/><m:endChr m:val="]"/><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-CA"/></w:rPr></m:ctrlPr></m:dPr><m:e><m:m><m:mPr><m:mcs><m:mc><m:mcPr><m:count m:val="4"/><m:mcJc m:val="center"/></m:mcPr></m:mc></m:mcs><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-CA"/></w:rPr></m:ctrlPr></m:mPr><m:mr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-CA"/></w:rPr><m:t>0</m:t></m:r></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-CA"/></w:rPr><m:t>7</m:t></m:r></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>8</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:e></m:mr><m:mr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-CA"/></w:rPr><m:t xml:space="preserve">0 </m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>5</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e></m:mr><m:mr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>0</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>2</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e></m:mr><m:mr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>0</m:t></m:r></m:e></m:mr></m:m></m:e></m:d><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:lang w:val="en-CA"/></w:rPr><w:br/></m:r></m:oMath></m:oMathPara>' },
    @{ Index = 5; Xml = '<m:oMathPara><m:oMath><m:d><m:dPr><m:begChr m:val="["/><m:endChr m:val="]"/><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-CA"/></w:rPr></m:ctrlPr></m:dPr><m:e><m:m><m:mPr><m:mcs><m:mc><m:mcPr><m:count m:val="4"/><m:mcJc m:val="center"/></m:mcPr></m:mc></m:mcs><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:lang w:val="en-CA"/></w:rPr></m:ctrlPr></m:mPr><m:mr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-CA"/></w:rPr><m:t>0</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-CA"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-CA"/></w:rPr><m:t>7</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-CA"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>12</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>8</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-CA"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-CA"/></w:rPr><m:t>0 </m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>5</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>7</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>0</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>2</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:e><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>∞</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:i/><w:lang w:val="de-DE"/></w:rPr></m:ctrlPr></m:e><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>0</m:t></m:r><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-CA"/></w:rPr></m:ctrlPr></m:e></m:mr></m:m><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-CA"/></w:rPr></m:ctrlPr></m:e></m:d></m:oMath></m:oMathPara>' }
)

foreach ($block in $mathBlocks) {
    $om = $d.OMaths.Item($block.Index)
    $om.Range.InsertXML($block.Xml)
}
